$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "Marco sala"
$ws.Range("B15").Value = "Lorenzo Canali | CGB Gamberoni"
$ws.Range("C15").Value = "Andrea  Roveda  | Pinguini Trentini"
$ws.Range("D15").Value = "Sebastiano Zoller | CGB Gamberoni"
$ws.Range("E15").Value = "Andrea Conzatti | FC SAVIGNANO"
$ws.Range("F15").Value = 'Lorenzo Mori` | Hellas Lazio'
